# 1st commit on 3rd July 2021
#
# - Rename sheet "USDBTC" -> "PAIR2" and replace its row-2 trade record.
# - Duplicate "PAIR2" into a new sheet "VNDUSD" (keeps identical headers,
#   column widths and cell styling) with two trade records of its own.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "PAIR2"

# --- PAIR2: overwrite the single data row (row 2) ---
$ws.Range("A2").Value = 7
$ws.Range("B2").Value = 6
$ws.Range("C2").Value = "14:58:01.136710"
$ws.Range("D2").Value = "Sell"
$ws.Range("E2").Value = "link 1"
$ws.Range("F2").Value = "link 2"
$ws.Range("G2").Value = "link 3"
$ws.Range("H2").Value = "link 4"
$ws.Range("I2").Value = "link 5"
$ws.Range("J2").Value = 3
$ws.Range("K2").Value = "this is my comment"
$ws.Range("L2").Value = 806
$ws.Range("M2").Value = 3

# --- Create VNDUSD as a copy of PAIR2 (preserves headers/col widths/styles) ---
$ws.Copy([System.Reflection.Missing]::Value, $ws)
$newws = $wb.Worksheets.Item(2)
$newws.Name = "VNDUSD"

# Row 2 on VNDUSD
$newws.Range("A2").Value = 7
$newws.Range("B2").Value = 6
$newws.Range("C2").Value = "14:57:30.000793"
$newws.Range("D2").Value = "Sell"
$newws.Range("E2").Value = "link 1"
$newws.Range("F2").Value = "link 2"
$newws.Range("G2").Value = "link 3"
$newws.Range("H2").Value = "link 4"
$newws.Range("I2").Value = "link 5"
$newws.Range("J2").Value = 3
$newws.Range("K2").Value = "this is my comment"
$newws.Range("L2").Value = 606
$newws.Range("M2").Value = 3

# Clone row 2's style down into row 3, then fill in its values
$newws.Range("A2:M2").Copy()
$newws.Range("A3:M3").PasteSpecial(-4122)  # xlPasteFormats

$newws.Range("A3").Value = 10
$newws.Range("B3").Value = 6
$newws.Range("C3").Value = "14:57:03.901108"
$newws.Range("D3").Value = "Sell"
$newws.Range("E3").Value = "link 1"
$newws.Range("F3").Value = "link 2"
$newws.Range("G3").Value = "link 3"
$newws.Range("H3").Value = "link 4"
$newws.Range("I3").Value = "link 5"
$newws.Range("J3").Value = 3
$newws.Range("K3").Value = "this is my comment"
$newws.Range("L3").Value = 406
$newws.Range("M3").Value = 3

# Keep PAIR2 as the selected/active sheet (matches original workbook state)
$ws.Select()
